$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 45
$ws.Range("H45").Value = 6503.778
$ws.Range("J45").Value = 9000
$ws.Range("L45").Value = 27000
$ws.Range("N45").Value = -27384
# Row 129
$ws.Range("H129").Value = 17106.064
$ws.Range("I129").Value = 1047.7333
$ws.Range("J129").Value = 22231.064
$ws.Range("K129").Value = 3143.199900000001
$ws.Range("L129").Value = 66693.192
$ws.Range("M129").Value = 1856.800099999999
$ws.Range("N129").Value = -76693.192
# Row 132
$ws.Range("H132").Value = 41815.895
$ws.Range("I132").Value = 31806.344
$ws.Range("K132").Value = 95419.03200000001
$ws.Range("M132").Value = -92889.03200000001
# Row 138
$ws.Range("H138").Value = 1996.9796
$ws.Range("I138").Value = 1193.25
$ws.Range("J138").Value = 2551.276
$ws.Range("K138").Value = 3579.75
$ws.Range("L138").Value = 7653.828
$ws.Range("M138").Value = 1560.25
$ws.Range("N138").Value = -17933.828

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13350.056
$ws.Range("I32").Value = 12550.83
$ws.Range("K32").Value = 12550.83
$ws.Range("M32").Value = -12263.83
# Row 88
$ws.Range("H88").Value = 12754649
$ws.Range("I88").Value = 66670132
$ws.Range("J88").Value = 2645496
$ws.Range("K88").Value = 66670132
$ws.Range("L88").Value = 2645496
$ws.Range("M88").Value = -66669726
$ws.Range("N88").Value = -2646308
# Row 91
$ws.Range("H91").Value = 12754649
$ws.Range("I91").Value = 66670132
$ws.Range("J91").Value = 2645496
$ws.Range("K91").Value = 66670132
$ws.Range("L91").Value = 2645496
$ws.Range("M91").Value = -66668728
$ws.Range("N91").Value = -2648304
# Row 97
$ws.Range("H97").Value = 1765.2
$ws.Range("I97").Value = 1628
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1628
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1132
$ws.Range("N97").Value = -3992

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2332.6667
$ws.Range("I94").Value = 2332.6667
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2332.6667
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1881.6667
$ws.Range("N94").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1750.55
$ws.Range("I31").Value = 1053
$ws.Range("J31").Value = 2235.288
$ws.Range("K31").Value = 1053
$ws.Range("L31").Value = 2235.288
$ws.Range("M31").Value = -758
$ws.Range("N31").Value = -2825.288
# Row 34
$ws.Range("H34").Value = 1750.55
$ws.Range("I34").Value = 1053
$ws.Range("J34").Value = 2235.288
$ws.Range("K34").Value = 1053
$ws.Range("L34").Value = 2235.288
$ws.Range("M34").Value = -851
$ws.Range("N34").Value = -2639.288
# Row 41
$ws.Range("H41").Value = 20646.334
$ws.Range("J41").Value = 28940
$ws.Range("L41").Value = 28940
$ws.Range("N41").Value = -29796
# Row 50
$ws.Range("H50").Value = 33281.5
$ws.Range("J50").Value = 38937.8
$ws.Range("L50").Value = 38937.8
$ws.Range("N50").Value = -40187.8
# Row 51
$ws.Range("H51").Value = 45490452
$ws.Range("J51").Value = 39498
$ws.Range("L51").Value = 39498
$ws.Range("N51").Value = -40970
# Row 59
$ws.Range("H59").Value = 31790
$ws.Range("J59").Value = 31790
$ws.Range("L59").Value = 31790
$ws.Range("N59").Value = -34080
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# Row 61
$ws.Range("H61").Value = 45490452
$ws.Range("J61").Value = 39498
$ws.Range("L61").Value = 39498
$ws.Range("N61").Value = -40194

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 10550.8
$ws.Range("I5").Value = 14644
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 43932
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -43820
$ws.Range("N5").Value = -3224
# Row 26
$ws.Range("H26").Value = 380.4
$ws.Range("J26").Value = 404
$ws.Range("L26").Value = 1212
$ws.Range("N26").Value = -1788
# Row 34
$ws.Range("H34").Value = 1568.9412
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 1798
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 5394
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -5562
# Row 68
$ws.Range("H68").Value = 1188.0121
$ws.Range("I68").Value = 802.96
$ws.Range("J68").Value = 1353.9828
$ws.Range("K68").Value = 2408.88
$ws.Range("L68").Value = 4061.9484
$ws.Range("M68").Value = -1597.88
$ws.Range("N68").Value = -5683.9484
# Row 71
$ws.Range("H71").Value = 1188.0121
$ws.Range("I71").Value = 802.96
$ws.Range("J71").Value = 1353.9828
$ws.Range("K71").Value = 7226.64
$ws.Range("L71").Value = 12185.8452
$ws.Range("M71").Value = -3170.64
$ws.Range("N71").Value = -20297.8452
# Row 86
$ws.Range("H86").Value = 800
$ws.Range("I86").Value = 733.3333
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2199.9999
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1013.9999
$ws.Range("N86").Value = -5372
# Row 89
$ws.Range("H89").Value = 800
$ws.Range("I89").Value = 733.3333
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6599.9997
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -671.9997000000003
$ws.Range("N89").Value = -20856
# Row 107
$ws.Range("H107").Value = 5332.7046
$ws.Range("I107").Value = 4037.2964
$ws.Range("J107").Value = 7390.1177
$ws.Range("K107").Value = 12111.8892
$ws.Range("L107").Value = 22170.3531
$ws.Range("M107").Value = -10191.8892
$ws.Range("N107").Value = -26010.3531
# Row 113
$ws.Range("H113").Value = 3466.9443
$ws.Range("I113").Value = 5817.316
$ws.Range("J113").Value = 840.05884
$ws.Range("K113").Value = 17451.948
$ws.Range("L113").Value = 2520.17652
$ws.Range("M113").Value = -15281.948
$ws.Range("N113").Value = -6860.17652
# Row 135
$ws.Range("H135").Value = 10550.8
$ws.Range("I135").Value = 14644
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 131796
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -129261
$ws.Range("N135").Value = -14070

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 27780172
$ws.Range("I132").Value = 41668170
$ws.Range("J132").Value = 4183.75
$ws.Range("K132").Value = 125004510
$ws.Range("L132").Value = 12551.25
$ws.Range("M132").Value = -125001980
$ws.Range("N132").Value = -17611.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1164.2858
$ws.Range("I93").Value = 825
$ws.Range("J93").Value = 1300
$ws.Range("K93").Value = 825
$ws.Range("L93").Value = 1300
$ws.Range("M93").Value = 423
$ws.Range("N93").Value = -3796
# Row 132
$ws.Range("H132").Value = 3547
$ws.Range("I132").Value = 3053.15
$ws.Range("J132").Value = 4066.842
$ws.Range("K132").Value = 9159.450000000001
$ws.Range("L132").Value = 12200.526
$ws.Range("M132").Value = -6629.450000000001
$ws.Range("N132").Value = -17260.526

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
